$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LBO Model")

# Sponsor Equity (% of Purchase Price) -> Sponsor Equity (% of Total Uses)
$ws.Range("A44").Value = "Sponsor Equity (% of Total Uses)"
# Sponsor Equity ($mm) formula: B13*B44 -> B27*B44
$ws.Range("B45").Formula = "=B27*B44"

# Senior Debt (% of Purchase Price) -> Senior Debt (% of Total Uses)
$ws.Range("A50").Value = "Senior Debt (% of Total Uses)"
# Senior Term Loan ($mm) formula: B13*B50 -> B27*B50
$ws.Range("B51").Formula = "=B27*B50"

# Subordinated Debt (% of Purchase Price) -> Subordinated Debt (% of Total Uses)
$ws.Range("A54").Value = "Subordinated Debt (% of Total Uses)"
# Subordinated Notes ($mm) formula: B13*B54 -> B27*B54
$ws.Range("B55").Formula = "=B27*B54"
